# Weekly price-sheet update.
# A brand-new weekly record is inserted as row 38 (pushing every existing
# record at/after the old row 38 down by one row), growing the used range
# from A1:R69 to A1:R70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 38; everything below (old rows 38-69)
# shifts down to rows 39-70, carrying its data/formatting along with it.
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with this week's record.
$ws.Range("A38").Value = 1
$ws.Range("B38").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C38").Value = "Arica y Parinacota"
$ws.Range("D38").Value = 44789
$ws.Range("E38").Value = 15
$ws.Range("F38").Value = 100112009
$ws.Range("G38").Value = "Acelga"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Segunda"
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 1500
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = 1500
$ws.Range("N38").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O38").Value = "Región de Arica y Parinacota"
$ws.Range("P38").Value = 500
$ws.Range("Q38").Value = 3
$ws.Range("R38").Value = "Hortaliza"
